$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header flags (F1 = FALSE, I1 = TRUE), mirrors B1/C1 boolean pattern ---
$ws.Range("F1").Value = $false
$ws.Range("I1").Value = $true

# --- Row 2 ("Knots" row): F2 = 1000, I2 = 500 ---
$ws.Range("F2").Value = 1000
$ws.Range("I2").Value = 500

# --- Rows 3-36: "F" (alt bias-corrected estimate), "G" (% diff vs column B, percent-styled,
#     color-scale conditional formatting), "I" (another estimate, black-font styled) ---
$ws.Cells.Item(3, 6).Value = 20664.4344417381
$ws.Cells.Item(3, 7).Formula = "=(F3-B3)/B3"
$ws.Cells.Item(3, 9).Value = 23367.3951
$ws.Cells.Item(4, 6).Value = 4132.52734004956
$ws.Cells.Item(4, 7).Formula = "=(F4-B4)/B4"
$ws.Cells.Item(4, 9).Value = 3596.54476
$ws.Cells.Item(5, 6).Value = 4132.52734004956
$ws.Cells.Item(5, 7).Formula = "=(F5-B5)/B5"
$ws.Cells.Item(5, 9).Value = 3596.54476
$ws.Cells.Item(6, 6).Value = 35737.5167355521
$ws.Cells.Item(6, 7).Formula = "=(F6-B6)/B6"
$ws.Cells.Item(6, 9).Value = 40961.9711
$ws.Cells.Item(7, 6).Value = 4132.52734004956
$ws.Cells.Item(7, 7).Formula = "=(F7-B7)/B7"
$ws.Cells.Item(7, 9).Value = 3596.54476
$ws.Cells.Item(8, 6).Value = 4132.52734004956
$ws.Cells.Item(8, 7).Formula = "=(F8-B8)/B8"
$ws.Cells.Item(8, 9).Value = 3596.54476
$ws.Cells.Item(9, 6).Value = 23366.18921461
$ws.Cells.Item(9, 7).Formula = "=(F9-B9)/B9"
$ws.Cells.Item(9, 9).Value = 27908.393
$ws.Cells.Item(10, 6).Value = 4132.52734004956
$ws.Cells.Item(10, 7).Formula = "=(F10-B10)/B10"
$ws.Cells.Item(10, 9).Value = 3596.54476
$ws.Cells.Item(11, 6).Value = 4132.52734004956
$ws.Cells.Item(11, 7).Formula = "=(F11-B11)/B11"
$ws.Cells.Item(11, 9).Value = 3596.54476
$ws.Cells.Item(12, 6).Value = 32124.372571352
$ws.Cells.Item(12, 7).Formula = "=(F12-B12)/B12"
$ws.Cells.Item(12, 9).Value = 36561.2308
$ws.Cells.Item(13, 6).Value = 4132.52734004956
$ws.Cells.Item(13, 7).Formula = "=(F13-B13)/B13"
$ws.Cells.Item(13, 9).Value = 3596.54476
$ws.Cells.Item(14, 6).Value = 4132.52734004956
$ws.Cells.Item(14, 7).Formula = "=(F14-B14)/B14"
$ws.Cells.Item(14, 9).Value = 3596.54476
$ws.Cells.Item(15, 6).Value = 32871.5709013146
$ws.Cells.Item(15, 7).Formula = "=(F15-B15)/B15"
$ws.Cells.Item(15, 9).Value = 37688.8637
$ws.Cells.Item(16, 6).Value = 4132.52734004956
$ws.Cells.Item(16, 7).Formula = "=(F16-B16)/B16"
$ws.Cells.Item(16, 9).Value = 3596.54476
$ws.Cells.Item(17, 6).Value = 4132.52734004956
$ws.Cells.Item(17, 7).Formula = "=(F17-B17)/B17"
$ws.Cells.Item(17, 9).Value = 3596.54476
$ws.Cells.Item(18, 6).Value = 30761.55811488
$ws.Cells.Item(18, 7).Formula = "=(F18-B18)/B18"
$ws.Cells.Item(18, 9).Value = 36031.7221
$ws.Cells.Item(19, 6).Value = 4132.52734004956
$ws.Cells.Item(19, 7).Formula = "=(F19-B19)/B19"
$ws.Cells.Item(19, 9).Value = 3596.54476
$ws.Cells.Item(20, 6).Value = 39919.3139589042
$ws.Cells.Item(20, 7).Formula = "=(F20-B20)/B20"
$ws.Cells.Item(20, 9).Value = 45993.1734
$ws.Cells.Item(21, 6).Value = 4132.52734004956
$ws.Cells.Item(21, 7).Formula = "=(F21-B21)/B21"
$ws.Cells.Item(21, 9).Value = 3596.54476
$ws.Cells.Item(22, 6).Value = 45065.064646497
$ws.Cells.Item(22, 7).Formula = "=(F22-B22)/B22"
$ws.Cells.Item(22, 9).Value = 50680.9009
$ws.Cells.Item(23, 6).Value = 4132.52734004956
$ws.Cells.Item(23, 7).Formula = "=(F23-B23)/B23"
$ws.Cells.Item(23, 9).Value = 3596.54476
$ws.Cells.Item(24, 6).Value = 66401.0771116699
$ws.Cells.Item(24, 7).Formula = "=(F24-B24)/B24"
$ws.Cells.Item(24, 9).Value = 76310.9947
$ws.Cells.Item(25, 6).Value = 4132.52734004956
$ws.Cells.Item(25, 7).Formula = "=(F25-B25)/B25"
$ws.Cells.Item(25, 9).Value = 3596.54476
$ws.Cells.Item(26, 6).Value = 42677.3287126539
$ws.Cells.Item(26, 7).Formula = "=(F26-B26)/B26"
$ws.Cells.Item(26, 9).Value = 50655.8537
$ws.Cells.Item(27, 6).Value = 4132.52734004956
$ws.Cells.Item(27, 7).Formula = "=(F27-B27)/B27"
$ws.Cells.Item(27, 9).Value = 3596.54476
$ws.Cells.Item(28, 6).Value = 33858.7089182581
$ws.Cells.Item(28, 7).Formula = "=(F28-B28)/B28"
$ws.Cells.Item(28, 9).Value = 42053.3878
$ws.Cells.Item(29, 6).Value = 4132.52734004956
$ws.Cells.Item(29, 7).Formula = "=(F29-B29)/B29"
$ws.Cells.Item(29, 9).Value = 3596.54476
$ws.Cells.Item(30, 6).Value = 38096.0617432645
$ws.Cells.Item(30, 7).Formula = "=(F30-B30)/B30"
$ws.Cells.Item(30, 9).Value = 44506.5614
$ws.Cells.Item(31, 6).Value = 4132.52734004956
$ws.Cells.Item(31, 7).Formula = "=(F31-B31)/B31"
$ws.Cells.Item(31, 9).Value = 3596.54476
$ws.Cells.Item(32, 6).Value = 54122.789507218
$ws.Cells.Item(32, 7).Formula = "=(F32-B32)/B32"
$ws.Cells.Item(32, 9).Value = 62920.1137
$ws.Cells.Item(33, 6).Value = 4132.52734004956
$ws.Cells.Item(33, 7).Formula = "=(F33-B33)/B33"
$ws.Cells.Item(33, 9).Value = 3596.54476
$ws.Cells.Item(34, 6).Value = 51720.6038295394
$ws.Cells.Item(34, 7).Formula = "=(F34-B34)/B34"
$ws.Cells.Item(34, 9).Value = 61673.5421
$ws.Cells.Item(35, 6).Value = 4132.52734004956
$ws.Cells.Item(35, 7).Formula = "=(F35-B35)/B35"
$ws.Cells.Item(35, 9).Value = 3596.54476
$ws.Cells.Item(36, 6).Value = 54009.0653919142
$ws.Cells.Item(36, 7).Formula = "=(F36-B36)/B36"
$ws.Cells.Item(36, 9).Value = 66091.1958

# --- Rows 37-38: new years 2018/2019 with only F filled in (B/C empty -> G is #DIV/0!) ---
$ws.Cells.Item(37, 1).Value = 2018
$ws.Cells.Item(37, 6).Value = 4132.52734004956
$ws.Cells.Item(37, 7).Formula = "=(F37-B37)/B37"
$ws.Cells.Item(38, 1).Value = 2019
$ws.Cells.Item(38, 6).Value = 76444.7337958478
$ws.Cells.Item(38, 7).Formula = "=(F38-B38)/B38"

# --- Formatting: G column picks up the existing "Percent" style (same as column E);
#     I column gets a black-font (explicit RGB) variant of the default font. ---
$percentFormat = $ws.Range("E3").NumberFormat
$ws.Range("G3:G36").NumberFormat = $percentFormat
$ws.Range("I3:I36").Font.Color = 0

# --- Conditional formatting: 3-color scale on the new G ratio column ---
$ws.Range("G3:G36").FormatConditions.AddColorScale(3) | Out-Null

# --- Leave the selection where the author left it when they saved ---
$ws.Range("O21").Select() | Out-Null
